$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C16").Value = "1143386178"
$ws.Range("D16").Value = "DANIELA DEL CARMEN MORA QUINTANA"
$ws.Range("E16").Value = "1707"
$ws.Range("F16").Value = 29509
$ws.Range("G16").Value = 900000
$ws.Range("C17").Value = "1143392521"
$ws.Range("D17").Value = "CRISTINA ISABEL ARRIETA AMELL"
$ws.Range("E17").Value = "1707"
$ws.Range("F17").Value = 29509
$ws.Range("G17").Value = 781242
$ws.Range("C18").Value = "1143386178"
$ws.Range("D18").Value = "DANIELA DEL CARMEN MORA QUINTANA"
$ws.Range("E18").Value = "1708"
$ws.Range("F18").Value = 29509
$ws.Range("G18").Value = 900000
$ws.Range("C19").Value = "1143392521"
$ws.Range("D19").Value = "CRISTINA ISABEL ARRIETA AMELL"
$ws.Range("E19").Value = "1708"
$ws.Range("F19").Value = 29509
$ws.Range("G19").Value = 781242
$ws.Range("C20").Value = "1143386178"
$ws.Range("D20").Value = "DANIELA DEL CARMEN MORA QUINTANA"
$ws.Range("E20").Value = "1709"
$ws.Range("F20").Value = 29509
$ws.Range("G20").Value = 900000
$ws.Range("C21").Value = "1143392521"
$ws.Range("D21").Value = "CRISTINA ISABEL ARRIETA AMELL"
$ws.Range("E21").Value = "1709"
$ws.Range("F21").Value = 29509
$ws.Range("G21").Value = 781242
$ws.Range("C22").Value = "1143386178"
$ws.Range("D22").Value = "DANIELA DEL CARMEN MORA QUINTANA"
$ws.Range("E22").Value = "1710"
$ws.Range("F22").Value = 29509
$ws.Range("G22").Value = 900000
$ws.Range("C23").Value = "1143392521"
$ws.Range("D23").Value = "CRISTINA ISABEL ARRIETA AMELL"
$ws.Range("E23").Value = "1710"
$ws.Range("F23").Value = 29509
$ws.Range("G23").Value = 781242
$ws.Range("C24").Value = "1143386178"
$ws.Range("D24").Value = "DANIELA DEL CARMEN MORA QUINTANA"
$ws.Range("E24").Value = "1711"
$ws.Range("F24").Value = 29509
$ws.Range("G24").Value = 900000
$ws.Range("C25").Value = "1143392521"
$ws.Range("D25").Value = "CRISTINA ISABEL ARRIETA AMELL"
$ws.Range("E25").Value = "1711"
$ws.Range("F25").Value = 29509
$ws.Range("G25").Value = 781242
$ws.Range("C26").Value = "1143386178"
$ws.Range("D26").Value = "DANIELA DEL CARMEN MORA QUINTANA"
$ws.Range("E26").Value = "1712"
$ws.Range("F26").Value = 29509
$ws.Range("G26").Value = 900000
$ws.Range("C27").Value = "1143392521"
$ws.Range("D27").Value = "CRISTINA ISABEL ARRIETA AMELL"
$ws.Range("E27").Value = "1712"
$ws.Range("F27").Value = 29509
$ws.Range("G27").Value = 781242
$ws.Range("C28").Value = "1043300081"
$ws.Range("D28").Value = "CAROLINA BARON RODRIGUEZ"
$ws.Range("E28").Value = "1810"
$ws.Range("F28").Value = 16666
$ws.Range("G28").Value = 781242
$ws.Range("C29").Value = "1043300081"
$ws.Range("D29").Value = "CAROLINA BARON RODRIGUEZ"
$ws.Range("E29").Value = "1811"
$ws.Range("F29").Value = 31249
$ws.Range("G29").Value = 781242
$ws.Range("C30").Value = "1043300081"
$ws.Range("D30").Value = "CAROLINA BARON RODRIGUEZ"
$ws.Range("E30").Value = "1812"
$ws.Range("F30").Value = 31249
$ws.Range("G30").Value = 781242
$ws.Range("C31").Value = "1043300081"
$ws.Range("D31").Value = "CAROLINA BARON RODRIGUEZ"
$ws.Range("E31").Value = "1901"
$ws.Range("F31").Value = 31249
$ws.Range("G31").Value = 781242
$ws.Range("C32").Value = "1043300081"
$ws.Range("D32").Value = "CAROLINA BARON RODRIGUEZ"
$ws.Range("E32").Value = "1902"
$ws.Range("F32").Value = 31249
$ws.Range("G32").Value = 781242
$ws.Range("C33").Value = "1043300081"
$ws.Range("D33").Value = "CAROLINA BARON RODRIGUEZ"
$ws.Range("E33").Value = "1903"
$ws.Range("F33").Value = 31249
$ws.Range("G33").Value = 781242
$ws.Range("C34").Value = "1043300081"
$ws.Range("D34").Value = "CAROLINA BARON RODRIGUEZ"
$ws.Range("E34").Value = "1904"
$ws.Range("F34").Value = 31249
$ws.Range("G34").Value = 781242
$ws.Range("C35").Value = "1043300081"
$ws.Range("D35").Value = "CAROLINA BARON RODRIGUEZ"
$ws.Range("E35").Value = "1905"
$ws.Range("F35").Value = 31249
$ws.Range("G35").Value = 781242
$ws.Range("C36").Value = "1043300081"
$ws.Range("D36").Value = "CAROLINA BARON RODRIGUEZ"
$ws.Range("E36").Value = "1906"
$ws.Range("F36").Value = 31249
$ws.Range("G36").Value = 781242
$ws.Range("C37").Value = "33109530"
$ws.Range("D37").Value = "LINET CECILIA TORRES ARROYO"
$ws.Range("E37").Value = "2102"
$ws.Range("F37").Value = 80400
$ws.Range("G37").Value = 2010000
$ws.Range("C38").Value = "73131582"
$ws.Range("D38").Value = "EDGAR YESID GONZALEZ SANTOS"
$ws.Range("E38").Value = "2102"
$ws.Range("F38").Value = 38000
$ws.Range("G38").Value = 950000
$ws.Range("C39").Value = "73132623"
$ws.Range("D39").Value = "OSCAR JOSE GUARDO SANTOYA"
$ws.Range("E39").Value = "2102"
$ws.Range("F39").Value = 228000
$ws.Range("G39").Value = 5700000
$ws.Range("C40").Value = "1143339688"
$ws.Range("D40").Value = "SARA ELENA MENDOZA PEREZ"
$ws.Range("E40").Value = "2102"
$ws.Range("F40").Value = 22533
$ws.Range("G40").Value = 1300000
$ws.Range("C41").Value = "1044926575"
$ws.Range("D41").Value = "SHIRLYS PAOLA MONTERO JIMENEZ"
$ws.Range("E41").Value = "2102"
$ws.Range("F41").Value = 35112
$ws.Range("G41").Value = 877803
$ws.Range("C42").Value = "22461412"
$ws.Range("D42").Value = "MARCELA PIA CARVAJALES GUTIERREZ"
$ws.Range("E42").Value = "2102"
$ws.Range("F42").Value = 375200
$ws.Range("G42").Value = 9380000
$ws.Range("C43").Value = "45553370"
$ws.Range("D43").Value = "LAURA ISABEL ARELLANO GARCIA"
$ws.Range("E43").Value = "2102"
$ws.Range("F43").Value = 104000
$ws.Range("G43").Value = 2600000
$ws.Range("C44").Value = "64582749"
$ws.Range("D44").Value = "OLGA LUCIA MONTERROZA PARRA"
$ws.Range("E44").Value = "2102"
$ws.Range("F44").Value = 35112
$ws.Range("G44").Value = 877803
$ws.Range("C45").Value = "1047432109"
$ws.Range("D45").Value = "YENIFER YESIS GAVIRIA REYES"
$ws.Range("E45").Value = "2102"
$ws.Range("F45").Value = 48000
$ws.Range("G45").Value = 1200000
$ws.Range("C46").Value = "1043339688"
$ws.Range("D46").Value = "SARA ELENA MENDOZA PEREZ"
$ws.Range("E46").Value = "2102"
$ws.Range("F46").Value = 13867
$ws.Range("G46").Value = 1300000
